$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 833474.2
$ws.Range("I6").Value = 1250078.8
$ws.Range("J6").Value = 265
$ws.Range("K6").Value = 3750236.4
$ws.Range("L6").Value = 795
$ws.Range("M6").Value = -3750124.4
$ws.Range("N6").Value = -1019
$ws.Range("H8").Value = 455.33334
$ws.Range("I8").Value = 3.625
$ws.Range("J8").Value = 971.5714
$ws.Range("K8").Value = 10.875
$ws.Range("L8").Value = 2914.7142
$ws.Range("M8").Value = 128.125
$ws.Range("N8").Value = -3192.7142
$ws.Range("H11").Value = 141.15
$ws.Range("I11").Value = 141.15
$ws.Range("K11").Value = 141.15
$ws.Range("M11").Value = -1.150000000000006
$ws.Range("H42").Value = 153.41667
$ws.Range("I42").Value = 14.333333
$ws.Range("J42").Value = 199.77777
$ws.Range("K42").Value = 42.999999
$ws.Range("L42").Value = 599.33331
$ws.Range("M42").Value = 187.000001
$ws.Range("N42").Value = -1059.33331
$ws.Range("H86").Value = 8297
$ws.Range("J86").Value = 8371.25
$ws.Range("L86").Value = 8371.25
$ws.Range("N86").Value = -10617.25
$ws.Range("H89").Value = 8297
$ws.Range("J89").Value = 8371.25
$ws.Range("L89").Value = 41856.25
$ws.Range("N89").Value = -53088.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 626.6667
$ws.Range("I4").Value = 190.125
$ws.Range("J4").Value = 1499.75
$ws.Range("K4").Value = 190.125
$ws.Range("L4").Value = 1499.75
$ws.Range("M4").Value = -74.125
$ws.Range("N4").Value = -1731.75
$ws.Range("H45").Value = 16028
$ws.Range("I45").Value = 26174
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 26174
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -25797
$ws.Range("N45").Value = -3254
$ws.Range("H80").Value = 49624.375
$ws.Range("J80").Value = 51999.285
$ws.Range("L80").Value = 51999.285
$ws.Range("N80").Value = -53995.285
$ws.Range("H83").Value = 49624.375
$ws.Range("J83").Value = 51999.285
$ws.Range("L83").Value = 155997.855
$ws.Range("N83").Value = -165981.855
$ws.Range("H106").Value = 36999.5
$ws.Range("J106").Value = 36999.5
$ws.Range("L106").Value = 36999.5
$ws.Range("N106").Value = -39523.5
$ws.Range("H132").Value = 3535.3
$ws.Range("I132").Value = 3200.353
$ws.Range("J132").Value = 5433.3335
$ws.Range("K132").Value = 9601.059000000001
$ws.Range("L132").Value = 16300.0005
$ws.Range("M132").Value = -7071.059000000001
$ws.Range("N132").Value = -21360.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 750
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 1000
$ws.Range("M8").Value = -860
$ws.Range("H134").Value = 3074.875
$ws.Range("I134").Value = 3074.875
$ws.Range("K134").Value = 9224.625
$ws.Range("M134").Value = -6689.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2526.111
$ws.Range("I58").Value = 1660.1538
$ws.Range("J58").Value = 4777.6
$ws.Range("K58").Value = 1660.1538
$ws.Range("L58").Value = 4777.6
$ws.Range("M58").Value = -1457.1538
$ws.Range("N58").Value = -5183.6
$ws.Range("H130").Value = 81999
$ws.Range("J130").Value = 81999
$ws.Range("L130").Value = 81999
$ws.Range("N130").Value = -92039
$ws.Range("H132").Value = 2398.4736
$ws.Range("I132").Value = 2364.2778
$ws.Range("J132").Value = 3014
$ws.Range("K132").Value = 7092.8334
$ws.Range("L132").Value = 9042
$ws.Range("M132").Value = -4562.8334
$ws.Range("N132").Value = -14102
$ws.Range("H134").Value = 7653.0444
$ws.Range("I134").Value = 4619.9736
$ws.Range("J134").Value = 24118.285
$ws.Range("K134").Value = 13859.9208
$ws.Range("L134").Value = 72354.855
$ws.Range("M134").Value = -11324.9208
$ws.Range("N134").Value = -77424.855
$ws.Range("H136").Value = 2526.111
$ws.Range("I136").Value = 1660.1538
$ws.Range("J136").Value = 4777.6
$ws.Range("K136").Value = 4980.4614
$ws.Range("L136").Value = 14332.8
$ws.Range("M136").Value = -2430.4614
$ws.Range("N136").Value = -19432.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 14.3
$ws.Range("I2").Value = 15.133333
$ws.Range("K2").Value = 90.799998
$ws.Range("M2").Value = 22.200002
$ws.Range("H7").Value = 1000104.6
$ws.Range("I7").Value = 113.71429
$ws.Range("K7").Value = 341.14287
$ws.Range("M7").Value = -229.14287
$ws.Range("H9").Value = 50052
$ws.Range("J9").Value = 83383.336
$ws.Range("L9").Value = 250150.008
$ws.Range("N9").Value = -250598.008
$ws.Range("H12").Value = 362.13794
$ws.Range("I12").Value = 290.22223
$ws.Range("J12").Value = 394.5
$ws.Range("K12").Value = 870.66669
$ws.Range("L12").Value = 1183.5
$ws.Range("M12").Value = -697.66669
$ws.Range("N12").Value = -1529.5
$ws.Range("H17").Value = 44.15
$ws.Range("J17").Value = 87
$ws.Range("L17").Value = 261
$ws.Range("N17").Value = -599
$ws.Range("H19").Value = 748.5
$ws.Range("I19").Value = 698.2
$ws.Range("K19").Value = 2094.6
$ws.Range("M19").Value = -1920.6
$ws.Range("H34").Value = 791.2222
$ws.Range("J34").Value = 3000
$ws.Range("L34").Value = 9000
$ws.Range("N34").Value = -9168
$ws.Range("H99").Value = 1598.6
$ws.Range("I99").Value = 1248.25
$ws.Range("K99").Value = 3744.75
$ws.Range("M99").Value = -1498.75
$ws.Range("H108").Value = 619.375
$ws.Range("I108").Value = 619.375
$ws.Range("K108").Value = 1858.125
$ws.Range("M108").Value = 1021.875
$ws.Range("H109").Value = 667666.7
$ws.Range("J109").Value = 3000
$ws.Range("L109").Value = 9000
$ws.Range("N109").Value = -11080
$ws.Range("H131").Value = 3362.3
$ws.Range("I131").Value = 1503
$ws.Range("K131").Value = 4509
$ws.Range("M131").Value = 531
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4024.3635
$ws.Range("I80").Value = 2645.3333
$ws.Range("J80").Value = 4541.5
$ws.Range("K80").Value = 2645.3333
$ws.Range("L80").Value = 4541.5
$ws.Range("M80").Value = -1647.3333
$ws.Range("N80").Value = -6537.5
$ws.Range("H83").Value = 4024.3635
$ws.Range("I83").Value = 2645.3333
$ws.Range("J83").Value = 4541.5
$ws.Range("K83").Value = 13226.6665
$ws.Range("L83").Value = 22707.5
$ws.Range("M83").Value = -8234.6665
$ws.Range("N83").Value = -32691.5
$ws.Range("H132").Value = 3917.3572
$ws.Range("I132").Value = 3236.9167
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 9710.750100000001
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -7180.750100000001
$ws.Range("N132").Value = -29060
$ws.Range("H139").Value = 73431.336
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20686
$ws.Range("H68").Value = 3499.923
$ws.Range("I68").Value = 3624.6667
$ws.Range("K68").Value = 3624.6667
$ws.Range("M68").Value = -2875.6667
$ws.Range("H71").Value = 3499.923
$ws.Range("I71").Value = 3624.6667
$ws.Range("K71").Value = 18123.3335
$ws.Range("M71").Value = -14379.3335
$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960
$ws.Range("H132").Value = 56198.79
$ws.Range("J132").Value = 5998.3335
$ws.Range("L132").Value = 17995.0005
$ws.Range("N132").Value = -23055.0005
$ws.Range("H136").Value = 5445.9287
$ws.Range("I136").Value = 5389.3
$ws.Range("K136").Value = 16167.9
$ws.Range("M136").Value = -13617.9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 27666.5
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802
$ws.Range("H130").Value = 35999.668
$ws.Range("J130").Value = 35999.668
$ws.Range("L130").Value = 35999.668
$ws.Range("N130").Value = -46039.668
$ws.Range("H131").Value = 155244
$ws.Range("J131").Value = 155244
$ws.Range("L131").Value = 155244
$ws.Range("N131").Value = -165324
$ws.Range("H132").Value = 2724.4138
$ws.Range("I132").Value = 2397.6538
$ws.Range("J132").Value = 5556.3335
$ws.Range("K132").Value = 7192.9614
$ws.Range("L132").Value = 16669.0005
$ws.Range("M132").Value = -4662.9614
$ws.Range("N132").Value = -21729.0005
